$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.039.76"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").Value = "3.150.15"
$ws.Range("E3").Value = "  +0.48%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.94%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "3.142.11"
$ws.Range("E8").Value = "  +0.27%  "

$ws.Range("E9").Value = "  -0.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.03%  "

$ws.Range("E13").Value = "  -2.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.27%  "

$ws.Range("D15").Value = "3.672.12"
$ws.Range("E15").Value = "  +0.51%  "

$ws.Range("E16").Value = "  -1.38%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.24%  "

$ws.Range("D18").Value = "63.907.75"
$ws.Range("E18").Value = "  +0.06%  "

$ws.Range("D19").Value = "3.148.89"
$ws.Range("E19").Value = "  +0.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.732"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.41%  "

$ws.Range("E23").Value = "  -0.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.72%  "

$ws.Range("E25").Value = "  +6.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.25%  "

$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.06%  "

$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("E31").Value = "  -0.46%  "

$ws.Range("E32").Value = "  +0.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.111"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.21%  "

$ws.Range("D35").Value = "0.0₃0839"
$ws.Range("E35").Value = "  -5.14%  "

$ws.Range("E36").Value = "  +1.16%  "

$ws.Range("E37").Value = "  -2.79%  "

$ws.Range("E38").Value = "  +0.09%  "

$ws.Range("E39").Value = "  -5.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "463.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.295"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.29%  "

$ws.Range("D44").Value = "2.926.93"
$ws.Range("E44").Value = "  +0.47%  "

$ws.Range("E45").Value = "  -0.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.93%  "

$ws.Range("E47").Value = "  -2.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.46%  "

$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("E50").Value = "  +2.17%  "

$ws.Range("E51").Value = "  -0.80%  "
